# Generate Report for Archive
#
# 1) Status text "Ready for handoff" -> "In Translation" everywhere it
#    appears (Overview!E2:F2 and the Status column, C2, on each per-locale
#    sheet).
# 2) Narrow the "Status" column(s) that held that text so the shorter
#    string still gets a tight autofit-style column width.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: zh-cn / de-de status cells (columns E & F, row 2) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value2 = "In Translation"
$overview.Range("F2").Value2 = "In Translation"
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

# --- Per-locale detail sheets: Status column (C2) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value2 = "In Translation"
$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value2 = "In Translation"
$dede.Columns.Item(3).ColumnWidth = 12.5
